$d = $word.ActiveDocument

# 1. Remove the "Meta description: ..." paragraph that follows the title heading.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description:*") {
        $p.Range.Delete()
        break
    }
}

# 2. Insert a new bold paragraph "Play Dragon Egg Free Slot Game - Review 2021"
#    right before the final paragraph (the one holding the DALLE prompt text).
$n = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($n)
$rng = $pLast.Range.Duplicate
$rng.Collapse(1)
$startPos = $rng.Start
$newHeading = "Play Dragon Egg Free Slot Game - Review 2021"
$rng.InsertBefore($newHeading + [char]13)

$headingRange = $d.Range($startPos, $startPos + $newHeading.Length)
$headingRange.Font.Bold = 1
$headingRange.Font.Italic = 0

# 3. Replace the DALLE image-prompt text in the final paragraph with the meta description text,
#    keeping the paragraph's existing (italic) formatting.
$oldPrompt = "Prompt for DALLE: Create a feature image for Dragon Egg, a slot game by Tom Horn, in a cartoon style. The image should feature a happy Maya warrior wearing glasses. The warrior should have a confident expression on their face and be holding a golden dragon egg in one hand, as if they have just won it in the slot game. The background should be a dark cave, with shadows of dragons visible in the background. The image should be eye-catching and convey the excitement of winning big in the game."
$newPrompt = "Read our review of Dragon Egg, a free slot game featuring dragons and symbols of wealth, with a straightforward interface and medium volatility."
$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newPrompt, 2)
